# Refresh the cryptocurrency market snapshot on Sheet1 (Price / Volume(1h)
# columns) to the latest values from the data feed, matching the commit
# "Updated cryptos list on Fri Sep 22 10:42:29 UTC 2023 with GitHub Actions".
#
# Every touched cell holds a plain string in the original workbook (e.g.
# "26.656.77", "  -0.47%  ") -- not a real number. Excel's COM layer will
# happily auto-coerce a string that merely *looks* like a plain decimal
# (e.g. "211.25") into a numeric cell when you assign .Value, which would
# silently change the cell's type/formatting (dropping trailing zeros such
# as "0.790" -> 0.79) and no longer match the source data. Forcing the
# cell's number format to Text ("@") before the assignment keeps every
# value a literal string, exactly like the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.656.77'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.597.15'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.25'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.245'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.68'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.820.83'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.605.49'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.99'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.646.62'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.67'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.91'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.40'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.97%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.28'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.663'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.92'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.297.46'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.16%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.842'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.91%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.790'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.74'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.734.15'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.15'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.62'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.50%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.51'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.77%  '
